$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D (which is empty) so columns E:I shift left to D:H
$ws.Columns("D").Delete()

# Update the selection to match the resulting state (entire column D selected)
$ws.Range("D1:D1048576").Select()
